# Modify protocol field in school transports' templates
#
# The "Αρ. Πρωτ.: ${protocol}" placeholder line gains a literal
# "Φ.15.1/" prefix right before the "${" that opens the merge field,
# so the rendered line reads "Αρ. Πρωτ.: Φ.15.1/${protocol}".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Αρ. Πρωτ.: `${",          # old text (literal "${" - no wildcards)
    $true,                     # MatchCase
    $false,                    # MatchWholeWord
    $false,                    # MatchWildcards
    $false,                    # MatchSoundsLike
    $false,                    # MatchAllWordForms
    $true,                     # Forward
    1,                         # Wrap (wdFindContinue)
    $false,                    # Format
    "Αρ. Πρωτ.: Φ.15.1/`${",   # new text
    2                          # Replace (wdReplaceAll)
)
